$d = $word.ActiveDocument

# The document has a handful of bold headings that start with a
# decorative emoji run (rendered in "Segoe UI Emoji") followed by a
# space and then the heading text in a separate run, e.g.:
#   [⭐][ BƯỚC 1 — Prompt ban đầu]
# The edit removes the emoji (and the space that glued it to the
# heading text), leaving just the heading text behind, e.g.:
#   [BƯỚC 1 — Prompt ban đầu]
# Do this for every occurrence of "⭐ " and "✅ " in the document.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("⭐ ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$find.Execute("✅ ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
